$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing label text (row 5, column A) ---
$ws.Range("A5").Value = "Implied minimum income (in `$/year)"

# --- Fix tiny floating point rounding on L5 ---
$ws.Range("L5").Value = 3291.56572730894

# --- Add new row 9: "Has not touched the sliders" ---
$ws.Range("A9").Value = "Has not touched the sliders"
$ws.Range("B9").Value = 0.408374634826883
$ws.Range("C9").Value = 0.388588909468523
$ws.Range("D9").Value = 0.399590058385227
$ws.Range("E9").Value = 0.384091647511391
$ws.Range("F9").Value = 0.410195834228784
$ws.Range("G9").Value = 0.302872499149178
$ws.Range("H9").Value = 0.371813108427663
$ws.Range("I9").Value = 0.437136028711529
$ws.Range("J9").Value = 0.303667681206425
$ws.Range("K9").Value = 0.475448989690047
$ws.Range("L9").Value = 0.388366284732463
$ws.Range("M9").Value = 0.406177169426307

# --- Add new row 10: "Touched sliders and satisfied" ---
$ws.Range("A10").Value = "Touched sliders and satisfied"
$ws.Range("B10").Value = 0.398118456901027
$ws.Range("C10").Value = 0.410488474534835
$ws.Range("D10").Value = 0.369081820440341
$ws.Range("E10").Value = 0.430038380653647
$ws.Range("F10").Value = 0.410781068299679
$ws.Range("G10").Value = 0.459008243270367
$ws.Range("H10").Value = 0.427013225872613
$ws.Range("I10").Value = 0.380331129442371
$ws.Range("J10").Value = 0.455290253425317
$ws.Range("K10").Value = 0.279967106328976
$ws.Range("L10").Value = 0.442866209980246
$ws.Range("M10").Value = 0.426810711448192
